$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-17 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-18 Sunday", 2) | Out-Null
$d.Content.Find.Execute("592×4=2368", $true, $false, $false, $false, $false, $true, 1, $false, "198×7=1386", 2) | Out-Null
$d.Content.Find.Execute("181×4=724", $true, $false, $false, $false, $false, $true, 1, $false, "360×2=720", 2) | Out-Null
$d.Content.Find.Execute("521×4=2084", $true, $false, $false, $false, $false, $true, 1, $false, "201×4=804", 2) | Out-Null
$d.Content.Find.Execute("557×8=4456", $true, $false, $false, $false, $false, $true, 1, $false, "892×5=4460", 2) | Out-Null
$d.Content.Find.Execute("364×6=2184", $true, $false, $false, $false, $false, $true, 1, $false, "304×5=1520", 2) | Out-Null
$d.Content.Find.Execute("854×7=5978", $true, $false, $false, $false, $false, $true, 1, $false, "392×5=1960", 2) | Out-Null
$d.Content.Find.Execute("967×2=1934", $true, $false, $false, $false, $false, $true, 1, $false, "104×7=728", 2) | Out-Null
$d.Content.Find.Execute("628×7=4396", $true, $false, $false, $false, $false, $true, 1, $false, "269×2=538", 2) | Out-Null
$d.Content.Find.Execute("620×3=1860", $true, $false, $false, $false, $false, $true, 1, $false, "852×2=1704", 2) | Out-Null
$d.Content.Find.Execute("130×4=520", $true, $false, $false, $false, $false, $true, 1, $false, "148×4=592", 2) | Out-Null
$d.Content.Find.Execute("502×3=1506", $true, $false, $false, $false, $false, $true, 1, $false, "507×7=3549", 2) | Out-Null
$d.Content.Find.Execute("288×8=2304", $true, $false, $false, $false, $false, $true, 1, $false, "107×4=428", 2) | Out-Null
$d.Content.Find.Execute("756×4=3024", $true, $false, $false, $false, $false, $true, 1, $false, "819×8=6552", 2) | Out-Null
$d.Content.Find.Execute("527×3=1581", $true, $false, $false, $false, $false, $true, 1, $false, "926×2=1852", 2) | Out-Null
$d.Content.Find.Execute("965×3=2895", $true, $false, $false, $false, $false, $true, 1, $false, "953×4=3812", 2) | Out-Null
$d.Content.Find.Execute("906×3=2718", $true, $false, $false, $false, $false, $true, 1, $false, "126×4=504", 2) | Out-Null
$d.Content.Find.Execute("130×7=910", $true, $false, $false, $false, $false, $true, 1, $false, "120×2=240", 2) | Out-Null
$d.Content.Find.Execute("962×5=4810", $true, $false, $false, $false, $false, $true, 1, $false, "245×5=1225", 2) | Out-Null
$d.Content.Find.Execute("174×5=870", $true, $false, $false, $false, $false, $true, 1, $false, "856×9=7704", 2) | Out-Null
$d.Content.Find.Execute("471×4=1884", $true, $false, $false, $false, $false, $true, 1, $false, "550×4=2200", 2) | Out-Null
$d.Content.Find.Execute("865×2=1730", $true, $false, $false, $false, $false, $true, 1, $false, "516×8=4128", 2) | Out-Null
$d.Content.Find.Execute("979×9=8811", $true, $false, $false, $false, $false, $true, 1, $false, "901×2=1802", 2) | Out-Null
$d.Content.Find.Execute("161×2=322", $true, $false, $false, $false, $false, $true, 1, $false, "672×2=1344", 2) | Out-Null
$d.Content.Find.Execute("569×3=1707", $true, $false, $false, $false, $false, $true, 1, $false, "280×5=1400", 2) | Out-Null
$d.Content.Find.Execute("473×8=3784", $true, $false, $false, $false, $false, $true, 1, $false, "389×9=3501", 2) | Out-Null
